$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.018.22'
$ws.Range('E2').Value = '  -0.71%  '
$ws.Range('D3').Value = '1.621.68'
$ws.Range('E3').Value = '  -1.47%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.31%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '214.98'
$ws.Range('E5').Value = '  -1.00%  '
$ws.Range('E6').Value = '  +0.21%  '
$ws.Range('E7').Value = '  +0.37%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.253'
$ws.Range('E8').Value = '  -1.13%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.0623'
$ws.Range('E9').Value = '  -0.74%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '20.11'
$ws.Range('E10').Value = '  +0.75%  '
$ws.Range('E11').Value = '  -0.15%  '
$ws.Range('D12').Value = '1.640.87'
$ws.Range('E12').Value = '  +0.15%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.11'
$ws.Range('E13').Value = '  -0.84%  '
$ws.Range('E14').Value = '  -0.12%  '
$ws.Range('B15').Value = 'Litecoin'
$ws.Range('C15').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '64.49'
$ws.Range('E15').Value = '  -4.72%  '
$ws.Range('B16').Value = 'WrappedBTC'
$ws.Range('C16').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D16').Value = '27.018.55'
$ws.Range('E16').Value = '  -0.63%  '
$ws.Range('D17').Value = '0.0₃0736'
$ws.Range('E17').Value = '  -0.34%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '215.86'
$ws.Range('E18').Value = '  -1.53%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.00'
$ws.Range('E19').Value = '  +0.21%  '
$ws.Range('E20').Value = '  +0.65%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.36'
$ws.Range('E21').Value = '  -1.19%  '
$ws.Range('E22').Value = '  -6.24%  '
$ws.Range('E23').Value = '  -2.56%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '147.53'
$ws.Range('E24').Value = '  -0.29%  '
$ws.Range('E25').Value = '  +0.34%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.26'
$ws.Range('E26').Value = '  -4.13%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.117'
$ws.Range('E27').Value = '  -0.73%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '15.54'
$ws.Range('E28').Value = '  -1.34%  '
$ws.Range('E29').Value = '  -0.95%  '
$ws.Range('E30').Value = '  -0.85%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.35'
$ws.Range('E31').Value = '  -1.19%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.98'
$ws.Range('E32').Value = '  -1.72%  '
$ws.Range('D33').Value = '1.331.40'
$ws.Range('E33').Value = '  +5.36%  '
$ws.Range('E34').Value = '  -1.10%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.46'
$ws.Range('E35').Value = '  +0.27%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0175'
$ws.Range('E36').Value = '  -1.26%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.543'
$ws.Range('E37').Value = '  -0.38%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.846'
$ws.Range('E38').Value = '  -0.22%  '
$ws.Range('E39').Value = '  +0.29%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.23'
$ws.Range('E40').Value = '  -0.91%  '
$ws.Range('E41').Value = '  -1.00%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '64.19'
$ws.Range('E42').Value = '  +3.69%  '
$ws.Range('D43').Value = '1.760.45'
$ws.Range('E43').Value = '  -1.48%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '5.22'
$ws.Range('E44').Value = '  -3.93%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '90.39'
$ws.Range('E45').Value = '  -1.38%  '
$ws.Range('E46').Value = '  -0.44%  '
$ws.Range('B47').Value = 'BabyDogeCoin'
$ws.Range('C47').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D47').Value = '0.0₆0107'
$ws.Range('E47').Value = '  -0.35%  '
$ws.Range('B48').Value = 'WEMIXToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.821'
$ws.Range('E48').Value = '  +22.52%  '
$ws.Range('E49').Value = '  -0.05%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0985'
$ws.Range('E50').Value = '  +1.24%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '7.56'
$ws.Range('E51').Value = '  -1.51%  '
